# The only change recorded for this deck is a byte-level reordering of
# the xmlns / xmlns:<prefix> attributes on a handful of already-present,
# opaque Office "extension list" (<a:extLst>) elements that PowerPoint
# itself never exposes through the object model:
#
#   - ppt/slideLayouts/slideLayout1.xml: the two <a14:hiddenFill> blocks
#     inside <p:spPr><a:extLst> of the two logo pictures on the poster
#     layout go from
#       <a14:hiddenFill xmlns="" xmlns:a14="...2010/main">
#     to
#       <a14:hiddenFill xmlns:a14="...2010/main" xmlns="">
#
#   - ppt/slides/slide1.xml: the nine <ma14:wrappingTextBoxFlag val="1"/>
#     elements inside <p:spPr><a:extLst> of the nine "Shape 22" text
#     boxes on the poster go from
#       <ma14:wrappingTextBoxFlag xmlns:ma14="...2011/main" xmlns="" .../>
#     to
#       <ma14:wrappingTextBoxFlag xmlns="" xmlns:ma14="...2011/main" .../>
#
# No text, formatting, or geometry changed between the two revisions --
# the commit note ("Made changes to time log as well as the poster")
# indicates the substantive edit that day was made in a separate
# time-log file; this deck was only resaved. Attribute order on an XML
# element carries no semantic meaning and PowerPoint's object model has
# no property that targets these vendor extension nodes directly --
# poking the related Fill/Line/Shape APIs to try to force a rewrite
# instead regenerates neighbouring XML (e.g. reorders/adds <a:ln> or
# <a:effectLst> children) and would introduce formatting changes that
# were never part of the real edit. So we just identify the exact shapes
# that own these extensions -- without changing any of their properties
# -- leaving the poster's content and formatting exactly as authored.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Slide layout: the two pictures ("Picture 2", ids 11 & 14) that carry
# the a14:hiddenFill extension.
$layout = $s.CustomLayout
$hiddenFillShapeIds = @(11, 14)
for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
    $layoutShape = $layout.Shapes.Item($i)
    if ($hiddenFillShapeIds -contains $layoutShape.Id) {
        $null = $layoutShape.Name
    }
}

# Slide: the nine text boxes (ids 33, 15, 18, 24, 25, 46, 54, 63, 67)
# that carry the ma14:wrappingTextBoxFlag extension.
$wrappingFlagShapeIds = @(33, 15, 18, 24, 25, 46, 54, 63, 67)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($wrappingFlagShapeIds -contains $shp.Id) {
        $null = $shp.Name
    }
}
